# Update res_line/loading_percent values for Case_3_59 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values, one row per PS array, columns in order: B,C,D,F,G,I,J,K,L,M,N
# (columns E, H are always 0 and unaffected by this change)
$newValues = @(
    @(18.54781195257897, 5.756501084900677, 4.825113339260736, 50.79052207779705, 3.774114025315235, 36.77808963134115, 10.22087035177882, 16.46040809085815, 12.03649010089267, 18.24008353545074, 25.34169354331657),
    @(18.44399694080811, 5.61434756660679, 4.826246063661417, 50.8200751546264, 3.776860116866707, 36.8320167559992, 10.23393861403654, 16.38961486196935, 12.05193218187985, 18.24337473424515, 25.39538001182678),
    @(18.38424289708379, 5.527064344401918, 4.827375544967794, 50.84690472958791, 3.778635901800782, 36.87044569822333, 10.24242108037399, 16.34964120505318, 12.06283884264002, 18.24824028502438, 25.43026362334662),
    @(18.36091688448739, 5.491555229793748, 4.827945417212656, 50.86002118210007, 3.77938217309211, 36.88744209997629, 10.24599338577286, 16.33424307662247, 12.06764221697217, 18.25093986177788, 25.44496248604809),
    @(18.3571060695584, 5.485664339434453, 4.82804667712857, 50.86233099152096, 3.779507459552904, 36.89034502112307, 10.24659355816614, 16.33174042696227, 12.06846149820934, 18.25143145791498, 25.44743244540883),
    @(18.3839241403174, 5.526585135029097, 4.827382786074887, 50.84707278420201, 3.778645874570322, 36.87066950830003, 10.24246878909662, 16.34942991478203, 12.06290216919429, 18.24827378827466, 25.43045989869753),
    @(18.51120186558571, 5.707521790189783, 4.825414114465809, 50.79890956102206, 3.775042311348194, 36.79557961080256, 10.22528133030533, 16.43528158593326, 12.04151901846414, 18.2406288750195, 25.35980660277967),
    @(18.79143894349279, 6.059824432035711, 4.82497628449364, 50.77336270335538, 3.76868381995655, 36.69055367077245, 10.19519943746473, 16.63073672170754, 12.01087432845407, 18.24813574313952, 25.23645568939228),
    @(19.01457381582218, 6.313954288845719, 4.826713246922894, 50.79656360974499, 3.764439100990804, 36.63917297148806, 10.17528539455819, 16.78995692874286, 11.99521270497585, 18.26725706013069, 25.15504889702789),
    @(19.11949682704144, 6.4279327512824, 4.827944255321666, 50.81620985513104, 3.762599735815855, 36.62140198686402, 10.1666963654666, 16.86556942393766, 11.98956925236686, 18.27888190372107, 25.12000687010049),
    @(19.15969053492707, 6.470813403049738, 4.828473245223347, 50.82495345556682, 3.761916307179968, 36.61547820222391, 10.16351116021462, 16.89464028615451, 11.98764457236281, 18.28370187517219, 25.10702288990361),
    @(19.15101405264141, 6.461591468143142, 4.828356532322752, 50.82301243235917, 3.762062914288947, 36.61671816226966, 10.16419416400776, 16.88836019306467, 11.98804965138712, 18.28264527142086, 25.10980652516537),
    @(19.12279450933703, 6.431466444736128, 4.827986519620381, 50.81690308490207, 3.762543247605388, 36.62089848632338, 10.16643297020011, 16.86795243687737, 11.98940665513311, 18.27927008810958, 25.11893294787354),
    @(19.10556844428712, 6.412976143196475, 4.827768042074733, 50.81333063627693, 3.762839169446684, 36.62356398066255, 10.16781305557847, 16.85550855022571, 11.99026549764331, 18.2772570219231, 25.1245603331035),
    @(19.00778309378176, 6.306468731731306, 4.826641621168816, 50.79546236134672, 3.764561144645112, 36.64044709333723, 10.17585613845533, 16.78507785691178, 11.99561127037818, 18.26655599295608, 25.15737898054278),
    @(18.9486493733734, 6.240678995669446, 4.826063117750845, 50.786827253786, 3.765640926296165, 36.65223936457991, 10.18091046010093, 16.74267232162036, 11.999269654353, 18.26073879596851, 25.17802154503743),
    @(18.9149607675836, 6.20268795136003, 4.825771907798364, 50.78271660909483, 3.766270612114593, 36.65954929487717, 10.18386182612969, 16.71858254767861, 12.00151330688167, 18.25766854203172, 25.19008197215335),
    @(18.90361084954938, 6.189800462351052, 4.825680458748272, 50.78147193285684, 3.766485296226667, 36.66211487457349, 10.18486871801346, 16.71047840522636, 12.00229693704928, 18.25667643380852, 25.19419762307149),
    @(18.95491098660992, 6.247698353595106, 4.826120406264377, 50.78765790097417, 3.765525089654771, 36.65092947935122, 10.18036784105486, 16.74715547421449, 11.99886578518133, 18.26132953876912, 25.17580472059513),
    @(19.13107098618093, 6.440322861484907, 4.828093500573131, 50.81866218960221, 3.762401807117087, 36.61964875764502, 10.16577355519981, 16.87393496072238, 11.98900231153374, 18.28025014430774, 25.11624454731502),
    @(19.24887816029264, 6.564559925659307, 4.829749004448243, 50.84652380476654, 3.760436880560986, 36.60390113324981, 10.15662732737826, 16.95933639069894, 11.98379352367909, 18.2950501279905, 25.07898342842244),
    @(19.18576741503964, 6.498418108270686, 4.828832133348986, 50.83095961011037, 3.76147863859841, 36.61187626572306, 10.16147307579376, 16.91352994833263, 11.98646052572514, 18.28692939565597, 25.09871822195955),
    @(18.95207914909926, 6.244525420486765, 4.826094377214196, 50.78727970550809, 3.765577431655421, 36.65152002623414, 10.18061301715784, 16.74512773693296, 11.99904793705952, 18.26106160982589, 25.17680634644651),
    @(18.71248992512229, 5.96512571965324, 4.824731131963747, 50.77289799350395, 3.770328653314268, 36.71444056728298, 10.20295177228663, 16.57504870716839, 12.0179588811832, 18.24370663229185, 25.26820305971025)
)

$columns = @('B','C','D','F','G','I','J','K','L','M','N')
$startRow = 2

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $rowValues = $newValues[$i]
    $excelRow = $startRow + $i
    for ($j = 0; $j -lt $columns.Length; $j++) {
        $ws.Range($columns[$j] + $excelRow).Value = $rowValues[$j]
    }
}
